$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I2").Value = "ArcGIS Platform"
$ws.Range("G2").Value = "Liity meihin"

$ws.Range("G2").Select()
